$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-18 Tuesday" "2025-11-19 Wednesday"
Replace-Text "70×75=" "14×60="
Replace-Text "41×59=" "51×45="
Replace-Text "61×20=" "65×45="
Replace-Text "69×90=" "64×60="
Replace-Text "71×38=" "77×19="
Replace-Text "95×58=" "45×15="
Replace-Text "11×83=" "52×47="
Replace-Text "88×73=" "96×51="
Replace-Text "46×36=" "50×76="
Replace-Text "93×19=" "30×81="
Replace-Text "78×21=" "48×68="
Replace-Text "23×69=" "50×68="
Replace-Text "51×60=" "22×43="
Replace-Text "59×40=" "84×22="
Replace-Text "88×97=" "14×50="
Replace-Text "62×92=" "17×63="
Replace-Text "85×71=" "94×80="
Replace-Text "26×63=" "53×84="
Replace-Text "51×71=" "90×86="
Replace-Text "28×25=" "98×38="
Replace-Text "20×54=" "22×59="
Replace-Text "91×90=" "44×93="
Replace-Text "16×53=" "29×34="
Replace-Text "22×89=" "98×32="
Replace-Text "93×88=" "46×87="
